$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.724.35'
$ws.Range('E2').Value = '  +5.63%  '
$ws.Range('D3').Value = '3.185.10'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '401.67'
$ws.Range('E5').Value = '  +3.30%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '109.41'
$ws.Range('E6').Value = '  +5.75%  '
$ws.Range('E7').Value = '  +1.44%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.618'
$ws.Range('E9').Value = '  +4.99%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '39.27'
$ws.Range('E10').Value = '  +5.78%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0890'
$ws.Range('E11').Value = '  +3.16%  '
$ws.Range('E12').Value = '  +1.69%  '
$ws.Range('D13').Value = '3.686.15'
$ws.Range('E13').Value = '  +2.54%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '19.08'
$ws.Range('E14').Value = '  +1.96%  '
$ws.Range('E15').Value = '  +2.04%  '
$ws.Range('E16').Value = '  +7.87%  '
$ws.Range('D17').Value = '3.188.51'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '10.57'
$ws.Range('E18').Value = '  -1.21%  '
$ws.Range('D19').Value = '54.618.85'
$ws.Range('E19').Value = '  +5.35%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.31'
$ws.Range('E20').Value = '  +3.30%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.0000102'
$ws.Range('E21').Value = '  +4.68%  '
$ws.Range('B22').Value = 'InternetComputer(DFINITY)'
$ws.Range('C22').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '12.94'
$ws.Range('E22').Value = '  +3.68%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '72.54'
$ws.Range('E23').Value = '  +3.53%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '276.25'
$ws.Range('E24').Value = '  +2.87%  '
$ws.Range('E25').Value = '  +4.54%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.00'
$ws.Range('E26').Value = '  -1.35%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.68'
$ws.Range('E27').Value = '  +6.47%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '27.88'
$ws.Range('E28').Value = '  +2.79%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.169'
$ws.Range('E29').Value = '  -0.54%  '
$ws.Range('B30').Value = 'Dai'
$ws.Range('C30').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.113'
$ws.Range('E31').Value = '  +2.92%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '11.07'
$ws.Range('E32').Value = '  +6.80%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0509'
$ws.Range('E33').Value = '  +13.38%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '36.95'
$ws.Range('E34').Value = '  +4.19%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.11'
$ws.Range('E35').Value = '  +1.55%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '51.34'
$ws.Range('E36').Value = '  +1.70%  '
$ws.Range('E37').Value = '  +6.03%  '
$ws.Range('E38').Value = '  -0.04%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.91'
$ws.Range('E39').Value = '  +12.05%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.11'
$ws.Range('E40').Value = '  +11.11%  '
$ws.Range('E41').Value = '  +3.32%  '
$ws.Range('E42').Value = '  +0.91%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '17.32'
$ws.Range('E43').Value = '  +2.48%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '131.25'
$ws.Range('E44').Value = '  +2.35%  '
$ws.Range('E45').Value = '  +1.43%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '22.25'
$ws.Range('E46').Value = '  -0.30%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.45'
$ws.Range('E47').Value = '  -2.21%  '
$ws.Range('E48').Value = '  -0.68%  '
$ws.Range('D49').Value = '2.095.33'
$ws.Range('E49').Value = '  +2.34%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0347'
$ws.Range('E50').Value = '  +9.15%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0508'
$ws.Range('E51').Value = '  +11.18%  '
